$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uppercase the existing student names
$ws.Range("A4").Value = "AARON TAN"
$ws.Range("A5").Value = "JOHN DOE"

# Add a new attendance-date column (D) with header "1-2"
$ws.Range("D3").Value = "1-2"
$ws.Range("D4").Value = 0.0
$ws.Range("D5").Value = 0.0
